# Weekly update for "Fruta / hortaliza" consolidated sheet.
# The rows in the data table (rows 2-30) get their date / volume / price
# columns (D, M, N, O, P, S) re-assigned according to a new weekly pairing,
# while every other column (market, product, category, unit, origin, ...)
# stays attached to the same physical row.
#
# Mapping: new row -> old row that supplies the D/M/N/O/P/S values.
$map = @{
    2  = 2
    3  = 4
    4  = 27
    5  = 9
    6  = 6
    7  = 23
    8  = 21
    9  = 13
    10 = 28
    11 = 16
    12 = 7
    13 = 12
    14 = 24
    15 = 25
    16 = 11
    17 = 18
    18 = 15
    19 = 17
    20 = 29
    21 = 3
    22 = 14
    23 = 20
    24 = 5
    25 = 10
    26 = 30
    27 = 26
    28 = 8
    29 = 19
    30 = 22
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "M", "N", "O", "P", "S")

# Snapshot the original values for the columns we are about to shuffle,
# before any writes happen (so reads never see already-mutated data).
$original = @{}
foreach ($row in $map.Values | Sort-Object -Unique) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $original[$row] = $rowVals
}

foreach ($newRow in ($map.Keys | Sort-Object)) {
    $oldRow = $map[$newRow]
    $src = $original[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $src[$col]
    }
}
